# Apply "Updated for design changes" revision to the PickAndPlace sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (and its tab) to reflect the new export date.
$ws.Name = "PickAndPlace_PCB1_2024-12-13"

# Row 2 — Designator "D1": updated Mid/Ref/Pad X & Y positions.
$ws.Range("D2").Value = "73.152mm"
$ws.Range("F2").Value = "73.152mm"
$ws.Range("E2").Value = "37.719mm"
$ws.Range("G2").Value = "37.719mm"
$ws.Range("I2").Value = "37.719mm"
$ws.Range("H2").Value = "71.979mm"

# Row 9 — Designator "RN1": updated Mid/Ref/Pad X & Y positions.
$ws.Range("D9").Value = "73.406mm"
$ws.Range("F9").Value = "73.406mm"
$ws.Range("E9").Value = "42.799mm"
$ws.Range("G9").Value = "42.799mm"
$ws.Range("H9").Value = "74.176mm"
$ws.Range("I9").Value = "41.449mm"
